$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-10 with the refreshed TPM-based NATMI values
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf16"
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.09476766666666665
$ws.Range("H2").Value = 0.284303
$ws.Range("I2").Value = 0.02336090049363864
$ws.Range("J2").Value = 0.02336090049363864
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.414593
$ws.Range("N2").Value = 4.243779
$ws.Range("O2").Value = 0.3478871232761722
$ws.Range("P2").Value = 0.3478871232761722
$ws.Range("Q2").Value = 0.134057677893
$ws.Range("R2").Value = 1.206519101037
$ws.Range("S2").Value = 0.008126956469872855
$ws.Range("T2").Value = 0.008126956469872857

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf16"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.09476766666666665
$ws.Range("H3").Value = 0.284303
$ws.Range("I3").Value = 0.02336090049363864
$ws.Range("J3").Value = 0.02336090049363864
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.337487333333333
$ws.Range("N3").Value = 7.012461999999999
$ws.Range("O3").Value = 0.5748520910875596
$ws.Range("P3").Value = 0.5748520910875596
$ws.Range("Q3").Value = 0.2215182204428889
$ws.Range("R3").Value = 1.993663983986
$ws.Range("S3").Value = 0.01342906249845657
$ws.Range("T3").Value = 0.01342906249845657

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf16"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.09476766666666665
$ws.Range("H4").Value = 0.284303
$ws.Range("I4").Value = 0.02336090049363864
$ws.Range("J4").Value = 0.02336090049363864
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.314161
$ws.Range("N4").Value = 0.942483
$ws.Range("O4").Value = 0.07726078563626818
$ws.Range("P4").Value = 0.07726078563626819
$ws.Range("Q4").Value = 0.02977230492766666
$ws.Range("R4").Value = 0.267950744349
$ws.Range("S4").Value = 0.001804881525309206
$ws.Range("T4").Value = 0.001804881525309206

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf16"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.187396333333333
$ws.Range("H5").Value = 6.562189
$ws.Range("I5").Value = 0.5392086761288135
$ws.Range("J5").Value = 0.5392086761288134
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.414593
$ws.Range("N5").Value = 4.243779
$ws.Range("O5").Value = 0.3478871232761722
$ws.Range("P5").Value = 0.3478871232761722
$ws.Range("Q5").Value = 3.094275541359
$ws.Range("R5").Value = 27.848479872231
$ws.Range("S5").Value = 0.1875837551840061
$ws.Range("T5").Value = 0.1875837551840061

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf16"
$ws.Range("C6").Value = "Fgfr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.187396333333333
$ws.Range("H6").Value = 6.562189
$ws.Range("I6").Value = 0.5392086761288135
$ws.Range("J6").Value = 0.5392086761288134
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.337487333333333
$ws.Range("N6").Value = 7.012461999999999
$ws.Range("O6").Value = 0.5748520910875596
$ws.Range("P6").Value = 0.5748520910875596
$ws.Range("Q6").Value = 5.113011222146445
$ws.Range("R6").Value = 46.017100999318
$ws.Range("S6").Value = 0.3099652350052031
$ws.Range("T6").Value = 0.3099652350052031

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf16"
$ws.Range("C7").Value = "Fgfr2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.187396333333333
$ws.Range("H7").Value = 6.562189
$ws.Range("I7").Value = 0.5392086761288135
$ws.Range("J7").Value = 0.5392086761288134
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.314161
$ws.Range("N7").Value = 0.942483
$ws.Range("O7").Value = 0.07726078563626818
$ws.Range("P7").Value = 0.07726078563626819
$ws.Range("Q7").Value = 0.6871946194763333
$ws.Range("R7").Value = 6.184751575287
$ws.Range("S7").Value = 0.04165968593960422
$ws.Range("T7").Value = 0.04165968593960422

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Fgf16"
$ws.Range("C8").Value = "Fgfr2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.774514666666667
$ws.Range("H8").Value = 5.323544
$ws.Range("I8").Value = 0.437430423377548
$ws.Range("J8").Value = 0.437430423377548
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.414593
$ws.Range("N8").Value = 4.243779
$ws.Range("O8").Value = 0.3478871232761722
$ws.Range("P8").Value = 0.3478871232761722
$ws.Range("Q8").Value = 2.510216025864
$ws.Range("R8").Value = 22.591944232776
$ws.Range("S8").Value = 0.1521764116222932
$ws.Range("T8").Value = 0.1521764116222932

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Fgf16"
$ws.Range("C9").Value = "Fgfr2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.774514666666667
$ws.Range("H9").Value = 5.323544
$ws.Range("I9").Value = 0.437430423377548
$ws.Range("J9").Value = 0.437430423377548
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.337487333333333
$ws.Range("N9").Value = 7.012461999999999
$ws.Range("O9").Value = 0.5748520910875596
$ws.Range("P9").Value = 0.5748520910875596
$ws.Range("Q9").Value = 4.147905556147555
$ws.Range("R9").Value = 37.331150005328
$ws.Range("S9").Value = 0.2514577935839
$ws.Range("T9").Value = 0.2514577935839

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Fgf16"
$ws.Range("C10").Value = "Fgfr2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.774514666666667
$ws.Range("H10").Value = 5.323544
$ws.Range("I10").Value = 0.437430423377548
$ws.Range("J10").Value = 0.437430423377548
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.314161
$ws.Range("N10").Value = 0.942483
$ws.Range("O10").Value = 0.07726078563626818
$ws.Range("P10").Value = 0.07726078563626819
$ws.Range("Q10").Value = 0.5574833021946666
$ws.Range("R10").Value = 5.017349719752
$ws.Range("S10").Value = 0.03379621817135477
$ws.Range("T10").Value = 0.03379621817135477

# Rows 11-13 (MuSCs/FAPs/Resolving-Mac combos) no longer exist in the refreshed data
$ws.Rows("11:13").Delete()
